# Arduino-MITCH Partslist edit: add a piezo buzzer (and related parts) to
# the PCB parts list, and drop the stray indent on the heading paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) The very first paragraph ("PCB Parts:") currently carries a
#    <w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr>. The target
#    drops that paragraph-properties block entirely.
# ---------------------------------------------------------------------
$headingXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t>PCB Parts:</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Paragraphs.Item(1).Range.InsertXML($headingXml)

# ---------------------------------------------------------------------
# 2) Append four new paragraphs at the very end of the body (after the
#    "Nand storage chip: KIOXIA TC58CVG2S0HRAIJ" bullet and before the
#    sectPr): a buzzer bullet, a magnetic-switch bullet, a blank
#    paragraph, and a "Screws:" paragraph.
# ---------------------------------------------------------------------
$newPartsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:color w:val="333333"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t xml:space="preserve">Buzzer: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:color w:val="333333"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t>PS1240P02BT</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:color w:val="333333"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t>F</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:color w:val="333333"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t xml:space="preserve">eatherweight </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:color w:val="333333"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t>Mag</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:color w:val="333333"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t>netic</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:color w:val="333333"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t xml:space="preserve"> Switch</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">Screws: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                <w:color w:val="212529"/>
                <w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/>
              </w:rPr>
              <w:t>&#188;-20</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$endOfDoc = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endOfDoc.InsertXML($newPartsXml)
